$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------
# Add two new year columns (K = 2021, L = 2022) mirroring the existing
# column J formatting, then fill in the reported values per row.
# -------------------------------------------------------------------

# Row 3 (thin bottom border spacer row, same style as existing J3)
$ws.Range("J3").Copy()
$ws.Range("K3:L3").PasteSpecial(-4122)

# Row 4 (header year row, same style as existing J4) with year values
$ws.Range("J4").Copy()
$ws.Range("K4:L4").PasteSpecial(-4122)
$ws.Range("K4").Value = 2021
$ws.Range("L4").Value = 2022

# Row 5 ("Всего" total row, same style as existing J5)
$ws.Range("J5").Copy()
$ws.Range("K5:L5").PasteSpecial(-4122)
$ws.Range("K5:L5").HorizontalAlignment = 1
$ws.Range("K5").Value = 272.60000000000002
$ws.Range("L5").Value = 292.19961890663211

# Row 6 ("в том числе:" blank sub-header row, same style as existing J6)
$ws.Range("J6").Copy()
$ws.Range("K6:L6").PasteSpecial(-4122)
$ws.Range("K6:L6").HorizontalAlignment = 1

# Row 7 (твердых), same style as existing J7
$ws.Range("J7").Copy()
$ws.Range("K7:L7").PasteSpecial(-4122)
$ws.Range("K7:L7").HorizontalAlignment = 1
$ws.Range("K7").Value = 98.1
$ws.Range("L7").Value = 99.522498012012946

# Row 8 (газообразных и жидких), same style as existing J8
$ws.Range("J8").Copy()
$ws.Range("K8:L8").PasteSpecial(-4122)
$ws.Range("K8:L8").HorizontalAlignment = 1
$ws.Range("K8").Value = 174.5
$ws.Range("L8").Value = 192.67712089461918

# Row 9 ("из них:" blank sub-header row), same style as existing J9
$ws.Range("J9").Copy()
$ws.Range("K9:L9").PasteSpecial(-4122)
$ws.Range("K9:L9").HorizontalAlignment = 1

# Row 10 (сернистого ангидрида), same style as existing J10
$ws.Range("J10").Copy()
$ws.Range("K10:L10").PasteSpecial(-4122)
$ws.Range("K10:L10").HorizontalAlignment = 1
$ws.Range("K10").Value = 75.599999999999994
$ws.Range("L10").Value = 88.011952928467494

# Row 11 (окиси углерода), same style as existing J11
$ws.Range("J11").Copy()
$ws.Range("K11:L11").PasteSpecial(-4122)
$ws.Range("K11:L11").HorizontalAlignment = 1
$ws.Range("K11").Value = 55.5
$ws.Range("L11").Value = 56.919430260413804

# Row 12 (окислов азота, bottom border row), same style as existing J12
$ws.Range("J12").Copy()
$ws.Range("K12:L12").PasteSpecial(-4122)
$ws.Range("K12:L12").HorizontalAlignment = 1
$ws.Range("K12").Value = 24.9
$ws.Range("L12").Value = 24.176373211436804

# Update the stored selection to match the authored workbook state.
$ws.Range("N5").Select()
